$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.065.57"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "3.421.80"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.29%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("E9").Value = "  +4.79%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.417"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("D12").Value = "4.004.37"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "3.418.30"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "62.108.93"
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "3.562.57"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "31.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "168.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "3.452.77"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("D47").Value = "2.548.29"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  -4.57%  "
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("E51").Value = "  -0.12%  "
